# Add two new GLM families ("Negative binomial" and "Beta") as new columns
# I and J in the status table header row, mirroring the style already used
# by the existing "Exponential" header in H2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells in row 2 (family names), columns I and J.
# H2 ("Exponential") uses the plain default style, so leave these as-is.
$ws.Range("I2").Value = "Negative binomial"
$ws.Range("J2").Value = "Beta"

# Widen the two new columns as in the target workbook (values chosen so the
# resulting OOXML column width, after this runtime's char/pixel rounding,
# lands as close as possible to the target widths of 16.58 and 5.32).
$ws.Columns.Item(9).ColumnWidth = 15.6
$ws.Columns.Item(10).ColumnWidth = 4.5
